$wb = $excel.ActiveWorkbook

# The Power Pivot worksheet-connection defined names lost their trailing "1"
# suffix (e.g. ...xlsxNode_Media1 -> ...xlsxNode_Media) when the template was
# refreshed/resaved; reproduce that rename here.
foreach ($definedName in $wb.Names) {
    if ($definedName.Name -match '^_xlcn\.WorksheetConnection_DraftTemplate\.xlsx.*1$') {
        $definedName.Name = $definedName.Name -replace '1$', ''
    }
}

# Add the new worksheet as the last tab (after RepeatingData) and rename it
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "DifferentTypes"

# Header row
$ws.Range("A1").Value = "Boolean"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "String"
$ws.Range("D1").Value = "Formula"

# Row 2
$ws.Range("A2").Value = $true
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "STRING"
$ws.Range("D2").Formula = "=B2+B3"

# Row 3
$ws.Range("A3").Value = $false
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "string"
$ws.Range("D3").Formula = "=C2+C3"

$ws.Range("A1:D3").Select() | Out-Null
